$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.68"
$ws.Range("E2").Value = "'-1.03%"
$ws.Range("D3").Value = "'37.03"
$ws.Range("E3").Value = "'-2.04%"
$ws.Range("D4").Value = "'5.122"
$ws.Range("E4").Value = "'-0.50%"
$ws.Range("D5").Value = "'0.07750"
$ws.Range("E5").Value = "'-2.13%"
$ws.Range("D6").Value = "'4.386"
$ws.Range("E6").Value = "'-0.55%"
$ws.Range("D7").Value = "'8.293"
$ws.Range("E7").Value = "'0.36%"
$ws.Range("D8").Value = "'1.858"
$ws.Range("E8").Value = "'-2.67%"
$ws.Range("E9").Value = "'-4.96%"
$ws.Range("D10").Value = "'0.9218"
$ws.Range("E10").Value = "'-0.64%"
$ws.Range("D11").Value = "'0.1136"
$ws.Range("E11").Value = "'-7.37%"
$ws.Range("D12").Value = "'0.1873"
$ws.Range("E12").Value = "'-2.44%"
$ws.Range("D13").Value = "'0.08823"
$ws.Range("E13").Value = "'-3.33%"
$ws.Range("D14").Value = "'0.03306"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("D15").Value = "'0.09552"
$ws.Range("E15").Value = "'-0.87%"
$ws.Range("D16").Value = "'0.001396"
$ws.Range("E16").Value = "'2.27%"
$ws.Range("D17").Value = "'0.006160"
$ws.Range("E17").Value = "'7.56%"
$ws.Range("D18").Value = "'3.386"
$ws.Range("E18").Value = "'-3.94%"
$ws.Range("D19").Value = "'0.3447"
$ws.Range("E19").Value = "'-0.12%"
$ws.Range("D20").Value = "'6.350"
$ws.Range("E20").Value = "'20.51%"
$ws.Range("D21").Value = "'0.1292"
$ws.Range("E21").Value = "'1.08%"
$ws.Range("D22").Value = "'0.2314"
$ws.Range("E22").Value = "'-10.75%"
$ws.Range("D23").Value = "'0.04348"
$ws.Range("E23").Value = "'-0.33%"
$ws.Range("E24").Value = "'-2.99%"
$ws.Range("D25").Value = "'0.004269"
$ws.Range("E25").Value = "'-0.73%"
$ws.Range("D26").Value = "'0.0001203"
$ws.Range("E26").Value = "'-1.33%"
$ws.Range("D27").Value = "'0.0002906"
$ws.Range("D39").Value = "'0.02118"
$ws.Range("E39").Value = "'-0.05%"
$ws.Range("D40").Value = "'0.04943"
$ws.Range("E40").Value = "'-4.97%"
$ws.Range("D41").Value = "'0.007596"
$ws.Range("E41").Value = "'0.14%"
$ws.Range("E42").Value = "'-0.63%"
$ws.Range("D43").Value = "'0.008581"
$ws.Range("E43").Value = "'-6.05%"
$ws.Range("D44").Value = "'0.002073"
$ws.Range("E44").Value = "'1.20%"
$ws.Range("D45").Value = "'0.008618"
$ws.Range("E45").Value = "'0.20%"
$ws.Range("D46").Value = "'0.00006579"
$ws.Range("E46").Value = "'-1.83%"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("D48").Value = "'0.003301"
$ws.Range("E48").Value = "'13.21%"
$ws.Range("D49").Value = "'0.001446"
$ws.Range("E49").Value = "'20.55%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.25%"
